$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.031.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.482.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.39"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.46"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.83%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.481.26"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.165"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.94"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.29%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.90"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.933.62"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.13%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.868.95"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.478.49"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.54"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "362.85"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.71%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.27%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.81"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -8.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0926"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.23%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "515.51"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.33%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.59%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.39%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.29"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.95"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.57"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.81%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -7.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.13"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.83"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.12%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.60"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0267"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.36%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.41%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.592"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.27%  "
